$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-08-06) was inserted before the existing
# row 19 (2021-07-26); every following row shifts down by one.
$ws.Rows.Item(19).Insert()

$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "Terminal La Palmera de La Serena"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44414
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112001
$ws.Range("G19").Value = "Berenjena"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 12500
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12750
$ws.Range("N19").Value = "$/caja 60 unidades"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 212
$ws.Range("Q19").Value = 60
$ws.Range("R19").Value = "Hortaliza"
